$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 323, pushing the existing rows 323-333 down to 324-334.
$ws.Rows(323).Insert()

# Populate the newly inserted row 323 with the new weekly price record.
$ws.Cells.Item(323, 1).Value = 11
$ws.Cells.Item(323, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(323, 3).Value = "Bíobío"
$ws.Cells.Item(323, 4).Value = 45147
$ws.Cells.Item(323, 5).Value = 8
$ws.Cells.Item(323, 6).Value = "Fruta"
$ws.Cells.Item(323, 7).Value = 100101
$ws.Cells.Item(323, 8).Value = "Berries"
$ws.Cells.Item(323, 9).Value = 100101007
$ws.Cells.Item(323, 10).Value = "Kiwi"
$ws.Cells.Item(323, 11).Value = "Hayward"
$ws.Cells.Item(323, 12).Value = "Especial"
$ws.Cells.Item(323, 13).Value = 220
$ws.Cells.Item(323, 14).Value = 14000
$ws.Cells.Item(323, 15).Value = 15000
$ws.Cells.Item(323, 16).Value = 14545
$ws.Cells.Item(323, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(323, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(323, 19).Value = 808
$ws.Cells.Item(323, 20).Value = 18
